# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 0.7527432677738641; E = 10.19245300693656;  G = 14.05633640148523 }
    3 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697;  G = 6.189590430959694 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697;  G = 6.189590430959694 }
    5 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697;  G = 6.189590430959694 }
    6 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 22.3905356188092;   E = 10.19245300693656;   G = 35.69412875252057 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Range("$col$row").Value = $values[$row][$col]
    }
}
